$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new text for row 10 (existing row, column C currently empty)
$ws.Range("C10").Value = "Different editors for different property types"

# Extend the numbered sequence in column B down to row 11
$ws.Range("B11").Formula = "=B10+1"

# Update the selected cell to match the new active cell C11
$ws.Range("C11").Select() | Out-Null
